$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17:B17").Copy()
$ws.Range("A18:B19").PasteSpecial(-4122)
$ws.Range("A18:B19").PasteSpecial(-4123)
$excel.CutCopyMode = 0

$ws.Range("A18").Value = 45727
$ws.Range("B18").Value = "Create a netflix like clone application with team"

$ws.Range("A19").Value = 45728
$ws.Range("B19").Value = "Covered Database concepts, "
$ws.Range("C19").Value = "Types of DB , NoSql,Cloud DB , In memory DB,ER relationships "

$ws.Range("B15").Select()
